# "refactor: change TVS diode"
#
# The BOM sheet had its "Quantity" and "Footprint" columns (C and D)
# swapped (Footprint now comes right after Designator, Quantity moved
# to column D), and the TVS diode row (designator D3) part number was
# updated from the old part (PESD5V2S2UT) to the new one (ESD5302).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap columns C and D in the header row (row 1): C1 becomes
#        "Footprint", D1 becomes "Quantity". Plain swap, same style. ---
$c1 = $ws.Cells.Item(1, 3).Value2
$d1 = $ws.Cells.Item(1, 4).Value2
$ws.Cells.Item(1, 3).Value = $d1
$ws.Cells.Item(1, 4).Value = $c1

# --- 2. Swap columns C (Quantity) and D (Footprint) for the data rows
#        2-20: column C becomes the (text) Footprint value, column D
#        becomes the (numeric) Quantity value. The leading "'" forces
#        the Footprint cell to keep a text/quote-prefixed style, just
#        like the rest of the text cells in the sheet. ---
for ($r = 2; $r -le 20; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value = "'" + $dVal
    $ws.Cells.Item($r, 4).Value = $cVal
}

# --- 3. Update the TVS diode (row 7 / designator D3) part number in
#        the Comment column from the old PESD5V2S2UT to ESD5302. ---
$ws.Range("A7").Value = "'ESD5302"

# --- 4. Column C used to be narrower (it held short quantity
#        numbers); now that it holds Footprint text it should match
#        the width of the other text columns. ---
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# --- 5. Reset the sheet view: drop the custom zoom level and the
#        lingering H12 selection, going back to a plain 100% view
#        with A1 selected. ---
$excel.ActiveWindow.Zoom = 100
$ws.Range("A1").Select()
